$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 731 (the panda post), shifting all rows below it up by one
$ws.Rows.Item(731).Delete()
